$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("git-manual")

$ws.Range("A13").Value = "git merge origin/<branchname>"
$ws.Range("B13").Value = "将分支<branchname>与当前分支合并"

$ws.Range("A14").Value = "git cat-file -p [<hash value>]"
$ws.Range("B14").Value = "查看.git/objects/<hash value前两位>/路径下保存的对象文件的内容"

$ws.Range("B6").Select()
